$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "id_scenario" column (column A) is unnecessary and is removed from
# the table. Deleting the entire column shifts every remaining column
# (id_region, id_building_type, id_building_component,
# id_building_component_option, id_action, unit, 2020, 2021) one position
# to the left.
$ws.Columns("A:A").Select() | Out-Null
$ws.Columns("A:A").Delete() | Out-Null
